# Adds a new forecast-creation-date column (AF, "2020-05-11") and a new
# target-date row (44, "2020-05-25") to both the "cases" and "deaths"
# sheets, plus fills in the newly-observed value for row 30 ("2020-05-11")
# in the "Observed" column (B) and the staircase of forecast values in the
# new AF column for rows 31-44.

$wb = $excel.ActiveWorkbook

# Per-sheet data: Observed (B30) value, and the AF column values for rows 31-44.
$sheetData = @{
    "cases" = @{
        B30 = 17939
        AF  = @(18988,19912,20805,21418,22188,22744,23266,23802,24401,25139,25733,26358,26842,27277)
    }
    "deaths" = @{
        B30 = 1770
        AF  = @(1886,1979,2084,2197,2296,2358,2418,2476,2564,2664,2758,2841,2888,2953)
    }
}

foreach ($sheetName in @("cases", "deaths")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $data = $sheetData[$sheetName]

    # --- New header cell AF1 = "2020-05-11" (reuse existing shared string; ---
    # --- must be forced to Text so Excel doesn't coerce it to a date.     ---
    $ws.Range("AF1").NumberFormat = "@"
    $ws.Range("AF1").Value = "2020-05-11"
    $ws.Range("AF1").Style = "Normal"

    # --- Rows 2-29: new column AF exists but stays empty. ---
    $ws.Range("AF2:AF29").NumberFormat = "General"
    $ws.Range("AF2:AF29").Style = "Normal"

    # --- Row 30: the "2020-05-11" target date is now observed. ---
    $ws.Range("B30").Value = $data.B30
    $ws.Range("AF30").NumberFormat = "General"
    $ws.Range("AF30").Style = "Normal"

    # --- Rows 31-44: new forecast values in column AF (the staircase). ---
    for ($i = 0; $i -lt $data.AF.Length; $i++) {
        $ws.Cells.Item(31 + $i, 32).Value = $data.AF[$i]
    }

    # --- New row 44 ("2020-05-25"): label + empty placeholders B..AE. ---
    $ws.Range("A44").NumberFormat = "@"
    $ws.Range("A44").Value = "2020-05-25"
    $ws.Range("A44").Style = "Normal"

    $ws.Range("B44:AE44").NumberFormat = "General"
    $ws.Range("B44:AE44").Style = "Normal"
}
